$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "kdfdj"
$ws.Range("D1").Value = "k"
$ws.Range("F1").Value = "j"
$ws.Range("G1").Value = "j"

$ws.Range("G1").Select()
